$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 10 - "Navigation - GPS Receiver": update cost, move/resize the
# specs table (it shrinks to make room for a new footnote box below)
# and add the footnote textbox itself.
# ---------------------------------------------------------------------
$s10 = $p.Slides.Item(10)
$tbl10 = $s10.Shapes.Item(2).Table

# Mark the "12,000" GPS receiver cost as already-owned hardware.
$tbl10.Cell(6, 3).Shape.TextFrame.TextRange.Text = "12,000*"

# Re-flow the row heights (table grows a touch to fit the new layout).
$tbl10.Rows.Item(1).Height = 68.46370078740158
$tbl10.Rows.Item(2).Height = 49.490787401574806
$tbl10.Rows.Item(3).Height = 34.23188976377953
$tbl10.Rows.Item(4).Height = 34.23188976377953
$tbl10.Rows.Item(5).Height = 68.46370078740158
$tbl10.Rows.Item(6).Height = 68.46370078740158

# Slide the whole table up so the footnote fits beneath it.
$s10.Shapes.Item(2).Top = 117.16354330708661

# New footnote textbox explaining the asterisk.
$note10 = $s10.Shapes.AddTextbox(1, 50.87488188976378, 467.34551181102364, 581.8909448818897, 29.081259842519685)
$note10.Name = "TextBox 2"
$note10.TextFrame.WordWrap = -1
$note10.TextFrame.AutoSize = 1
$note10.Fill.Visible = 0
$note10.TextFrame.TextRange.Text = "*Currently owned and wouldn" + [char]0x2019 + "t contribute to total budget"

# ---------------------------------------------------------------------
# Slide 11 - "Navigation - GPS Antenna": the antenna (NovaTel antenna)
# is already owned, so its cost cell changes from N/A to "Already owned".
# ---------------------------------------------------------------------
$s11 = $p.Slides.Item(11)
$tbl11 = $s11.Shapes.Item(2).Table
$tbl11.Cell(6, 4).Shape.TextFrame.TextRange.Text = "Already owned"
